# Insert a new data row above row 279, pushing existing rows 279:313 down to 280:314
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(279).Insert()

# Populate the newly inserted row 279 with the new weekly record
$ws.Range("A279").Value = 11
$ws.Range("B279").Value = "Vega Monumental Concepción"
$ws.Range("C279").Value = "Bíobío"
$ws.Range("D279").Value = 45142
$ws.Range("E279").Value = 8
$ws.Range("F279").Value = 100112003
$ws.Range("G279").Value = "Ajo"
$ws.Range("H279").Value = "Chino"
$ws.Range("I279").Value = "Primera"
$ws.Range("J279").Value = 240
$ws.Range("K279").Value = 18000
$ws.Range("L279").Value = 19000
$ws.Range("M279").Value = 18500
$ws.Range("N279").Value = "$/caja 10 kilos"
$ws.Range("O279").Value = "China"
$ws.Range("P279").Value = 1850
$ws.Range("Q279").Value = 10
$ws.Range("R279").Value = "Hortaliza"
